$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Elementos del Proyecto")

# New "Interfaz de Usuario" rows (J/K/L, rows 10-18) finishing the prototype
# documentation traceability table ("todos los doc terminados").
$rows = @(
    @{ Row = 10; Nombre = "Prototipo Login";            Id = 7;  Casos = "11,12" },
    @{ Row = 11; Nombre = "prototipo Mi cuenta";         Id = 8;  Casos = 22 },
    @{ Row = 12; Nombre = "prototipo agregar mascota";   Id = 9;  Casos = 17 },
    @{ Row = 13; Nombre = "prototipo buscar";            Id = 10; Casos = 19 },
    @{ Row = 14; Nombre = "prototipo evento";            Id = 11; Casos = 15 },
    @{ Row = 15; Nombre = "prototipo filtros";           Id = 12; Casos = 20 },
    @{ Row = 16; Nombre = "prototipo mis mascotas";      Id = 13; Casos = "16,18" },
    @{ Row = 17; Nombre = "prototipo nueva cuenta";      Id = 14; Casos = 10 },
    @{ Row = 18; Nombre = "prototipo servicios";         Id = 15; Casos = "13,14" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 10).Value = $r.Nombre
    $ws.Cells.Item($r.Row, 11).Value = $r.Id
    $ws.Cells.Item($r.Row, 12).Value = $r.Casos
}

# Update the sheet's view / selection to match where the author left off.
$ws.Activate()
$ws.Range("K19").Select()
